$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.132280349731445
$ws.Range("B1").Value = 2.558017730712891
$ws.Range("C1").Value = 6.153486728668213
$ws.Range("D1").Value = 2.153107881546021
$ws.Range("E1").Value = 1.240064024925232
